$wb = $excel.ActiveWorkbook

$wsLibrary      = $wb.Worksheets.Item("Library")
$wsFormula      = $wb.Worksheets.Item("Library_Formula")
$wsLabels       = $wb.Worksheets.Item("Formula Libraries Labels")

# --- Add the new "backTesting" row to Library_Formula (row 4) ---
$newRow = $wsFormula.Range("A4:G4")
$newRow.Font.Name = "Trebuchet MS"
$newRow.Font.Size = 10
$newRow.Font.ColorIndex = 1

$wsFormula.Cells.Item(4, 1).Value = "CREATE/MODIFY"
$wsFormula.Cells.Item(4, 2).Value = "LIB_RISK_EW"
$wsFormula.Cells.Item(4, 3).Value = "backTesting"
$wsFormula.Cells.Item(4, 5).Value = "Date"

# --- Update the selection (active cell) on every sheet, finishing on ---
# --- Library_Formula so it remains the active tab.                    ---
$wsLibrary.Range("B12").Select() | Out-Null
$wsLabels.Range("B22").Select() | Out-Null
$wsFormula.Range("C5").Select() | Out-Null
